# Edit script: rewrite "Light--A Cosmic Messenger" essay into a chemistry
# essay ("The Enchanting Realm of Chemistry..."), update author/e-mail,
# rewrite the body + summary text, and append a trailing blank paragraph.

$d = $word.ActiveDocument
$brk = [char]11   # manual line break (w:br) inside a paragraph

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Light--A Cosmic Messenger", $false, $false, $false, $false, $false, $true, 1, $false, "The Enchanting Realm of Chemistry: Unraveling the Secrets of Matter", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Author name paragraph -> "Dr. Sophia Alexander"
# ---------------------------------------------------------------------
$authorPara = $d.Paragraphs.Item(2)
$authorStart = $authorPara.Range.Start
$authorEnd = $authorPara.Range.End
$authorRange = $d.Range($authorStart, $authorEnd - 1)
$authorRange.Text = "Dr. Sophia Alexander"

# ---------------------------------------------------------------------
# 3. E-mail paragraph -> "sophiaa@educonnect.org"
# ---------------------------------------------------------------------
$emailPara = $d.Paragraphs.Item(3)
$emailStart = $emailPara.Range.Start
$emailEnd = $emailPara.Range.End
$emailRange = $d.Range($emailStart, $emailEnd - 1)
$emailRange.Text = "sophiaa@educonnect.org"

# ---------------------------------------------------------------------
# 4. Big intro/body paragraph (paragraph 5)
# ---------------------------------------------------------------------
$bodyText = 'In the vast expanse of scientific disciplines, chemistry stands out as a captivating realm that unveils the intricate world of matter and its transformations. It is a subject that intertwines the tangible and the abstract, encompassing the study of substances, their properties, and the intricate dance of chemical reactions that shape our world.'
$bodyText = $bodyText + $brk + $brk + 'Chemistry is the science that delves into the fundamental building blocks of matter, exploring the structure and behavior of atoms, molecules, and compounds. It''s a discipline that seeks to unravel the enigmas of chemical reactions, unmasking the intricate mechanisms that govern how substances interact and transform. Through the lens of chemistry, we gain invaluable insights into the composition of the universe, the intricate workings of living organisms, and the countless chemical processes that shape our everyday lives.'
$bodyText = $bodyText + $brk + $brk + 'The study of chemistry opens doors to a world of wonder and practical applications. From understanding the intricate interactions of atoms to harnessing the power of chemical reactions for life-saving medicines and groundbreaking technologies, chemistry plays a pivotal role in shaping our future. It permeates every aspect of our lives, from the food we eat to the clothes we wear and the medicines that heal us, making it an indispensable tool for navigating the complexities of the modern world.'
$bodyText = $bodyText + $brk + $brk + 'Introduction Continued:'
$bodyText = $bodyText + $brk + $brk + 'Chemistry is a subject that captivates the imagination and encourages a spirit of exploration. It invites us to question the world around us, to probe the mysteries of matter, and to uncover the underlying principles that govern chemical phenomena. Through hands-on experiments and thought-provoking discussions, chemistry classes provide students with an opportunity to engage their curiosity, develop critical thinking skills, and cultivate a deep appreciation for the natural world.'
$bodyText = $bodyText + $brk + $brk + 'From the alchemists of ancient times to the pioneering chemists of today, the pursuit of chemical knowledge has driven countless individuals to push the boundaries of human understanding. The discoveries made in chemistry have transformed our world, leading to advancements in medicine, agriculture, materials science, and countless other fields. It is a discipline that has played a pivotal role in shaping human history and continues to hold immense promise for addressing the challenges of the future.'
$bodyText = $bodyText + $brk + $brk + 'Introduction Continued:'
$bodyText = $bodyText + $brk + $brk + 'The study of chemistry is not without its challenges. It demands a willingness to grapple with abstract concepts, to persevere through difficult problems, and to embrace the complexities of a subject that is constantly evolving. However, the rewards of this endeavor are immense. Chemistry opens doors to countless career opportunities, from research and development to medicine, engineering, and environmental science. More importantly, it equips individuals with a deeper understanding of the world around them and empowers them to make informed decisions about the future.'

$bodyPara = $d.Paragraphs.Item(5)
$bodyStart = $bodyPara.Range.Start
$bodyEnd = $bodyPara.Range.End
$bodyRange = $d.Range($bodyStart, $bodyEnd - 1)
$bodyRange.Text = $bodyText

# ---------------------------------------------------------------------
# 5. Summary body paragraph (last paragraph, after "Summary" heading)
# ---------------------------------------------------------------------
$summaryText = 'Chemistry is a captivating and challenging subject that delves into the intricate world of matter and its transformations. It is a field of study that encompasses the study of substances, their properties, and the intricate dance of chemical reactions that shape our world. The study of chemistry opens doors to a world of wonder and practical applications, from understanding the composition of the universe to harnessing the power of chemical reactions for life-saving medicines and groundbreaking technologies. Chemistry is a subject that captivates the imagination, encourages exploration, and plays a pivotal role in shaping our future.'

$summaryPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$summaryStart = $summaryPara.Range.Start
$summaryEnd = $summaryPara.Range.End
$summaryRange = $d.Range($summaryStart, $summaryEnd - 1)
$summaryRange.Text = $summaryText

# ---------------------------------------------------------------------
# 6. Append a trailing blank paragraph after the summary
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 7. Normalize the font everywhere: TimesNewToman -> Times New Roman
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $d.Paragraphs.Item($i).Range.Font.Name = "Times New Roman"
}
